# Update the "想去人数" (F column) figures on the "展览" and "全部类型"
# sheets to reflect the latest scrape (commit 456a3b4):
#   - Row 9  (合肥·首届AT次元时代动漫游戏嘉年华): F9  3723 -> 3727
#   - Row 10 (合肥·W·A第五人格同人only2.0):        F10   66 -> 67

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F9").Value = 3727
    $ws.Range("F10").Value = 67
}
